$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 1834
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1834
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1834
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = -2486

# Row 123
$ws.Range("H123").Value = 40113.332
$ws.Range("J123").Value = 40113.332
$ws.Range("L123").Value = 40113.332
$ws.Range("N123").Value = -49913.332

# Row 137
$ws.Range("H137").Value = 17393.354
$ws.Range("I137").Value = 1152.8182
$ws.Range("J137").Value = 47167.668
$ws.Range("K137").Value = 3458.4546
$ws.Range("L137").Value = 141503.004
$ws.Range("M137").Value = -908.4546
$ws.Range("N137").Value = -146603.004

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1511.4286
$ws.Range("I2").Value = 1176
$ws.Range("J2").Value = 2350
$ws.Range("K2").Value = 1176
$ws.Range("L2").Value = 2350
$ws.Range("M2").Value = -1063
$ws.Range("N2").Value = -2576

# Row 32
$ws.Range("H32").Value = 12335.37
$ws.Range("I32").Value = 11954.919
$ws.Range("K32").Value = 11954.919
$ws.Range("M32").Value = -11667.919

# Row 45
$ws.Range("H45").Value = 1895473.4
$ws.Range("I45").Value = 2393840
$ws.Range("J45").Value = 1680
$ws.Range("K45").Value = 2393840
$ws.Range("L45").Value = 1680
$ws.Range("M45").Value = -2393463
$ws.Range("N45").Value = -2434

# Row 116
$ws.Range("H116").Value = 1511.4286
$ws.Range("I116").Value = 1176
$ws.Range("J116").Value = 2350
$ws.Range("K116").Value = 1176
$ws.Range("L116").Value = 2350
$ws.Range("M116").Value = 1118
$ws.Range("N116").Value = -6938

# Row 132
$ws.Range("H132").Value = 1890.1666
$ws.Range("I132").Value = 1473.9706
$ws.Range("J132").Value = 2900.9285
$ws.Range("K132").Value = 4421.9118
$ws.Range("L132").Value = 8702.7855
$ws.Range("M132").Value = -1891.9118
$ws.Range("N132").Value = -13762.7855

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1511.4286
$ws.Range("I3").Value = 1176
$ws.Range("J3").Value = 2350
$ws.Range("K3").Value = 1176
$ws.Range("L3").Value = 2350
$ws.Range("M3").Value = -1062
$ws.Range("N3").Value = -2578

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 14207.75
$ws.Range("J16").Value = 3660.6667
$ws.Range("L16").Value = 3660.6667
$ws.Range("N16").Value = -4234.6667

# Row 31
$ws.Range("H31").Value = 5633.8945
$ws.Range("I31").Value = 1839
$ws.Range("J31").Value = 8393.817999999999
$ws.Range("K31").Value = 1839
$ws.Range("L31").Value = 8393.817999999999
$ws.Range("M31").Value = -1544
$ws.Range("N31").Value = -8983.817999999999

# Row 34
$ws.Range("H34").Value = 5633.8945
$ws.Range("I34").Value = 1839
$ws.Range("J34").Value = 8393.817999999999
$ws.Range("K34").Value = 1839
$ws.Range("L34").Value = 8393.817999999999
$ws.Range("M34").Value = -1637
$ws.Range("N34").Value = -8797.817999999999

# Row 99
$ws.Range("H99").Value = 2607.2307
$ws.Range("I99").Value = 1460
$ws.Range("J99").Value = 3324.25
$ws.Range("K99").Value = 1460
$ws.Range("L99").Value = 3324.25
$ws.Range("M99").Value = 38
$ws.Range("N99").Value = -6320.25

# Row 113
$ws.Range("H113").Value = 14207.75
$ws.Range("J113").Value = 3660.6667
$ws.Range("L113").Value = 3660.6667
$ws.Range("N113").Value = -8000.6667

# Row 126
$ws.Range("H126").Value = 2607.2307
$ws.Range("I126").Value = 1460
$ws.Range("J126").Value = 3324.25
$ws.Range("K126").Value = 4380
$ws.Range("L126").Value = 9972.75
$ws.Range("M126").Value = -1910
$ws.Range("N126").Value = -14912.75

# Row 132
$ws.Range("H132").Value = 3667.56
$ws.Range("I132").Value = 1621.8158
$ws.Range("J132").Value = 10145.75
$ws.Range("K132").Value = 4865.4474
$ws.Range("L132").Value = 30437.25
$ws.Range("M132").Value = -2335.4474
$ws.Range("N132").Value = -35497.25

# Row 134
$ws.Range("H134").Value = 1324986.4
$ws.Range("I134").Value = 2314.4
$ws.Range("J134").Value = 7938346
$ws.Range("K134").Value = 6943.200000000001
$ws.Range("L134").Value = 23815038
$ws.Range("M134").Value = -4408.200000000001
$ws.Range("N134").Value = -23820108

$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 1442.3235
$ws.Range("I132").Value = 1201
$ws.Range("J132").Value = 1474.5
$ws.Range("K132").Value = 10809
$ws.Range("L132").Value = 13270.5
$ws.Range("M132").Value = -8279
$ws.Range("N132").Value = -18330.5

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2874.5454
$ws.Range("I122").Value = 2808.3572
$ws.Range("J122").Value = 2990.375
$ws.Range("K122").Value = 8425.071599999999
$ws.Range("L122").Value = 8971.125
$ws.Range("M122").Value = -5975.071599999999
$ws.Range("N122").Value = -13871.125

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7078.8
$ws.Range("I7").Value = 4916.625
$ws.Range("J7").Value = 9549.857
$ws.Range("K7").Value = 4916.625
$ws.Range("L7").Value = 9549.857
$ws.Range("M7").Value = -4804.625
$ws.Range("N7").Value = -9773.857

# Row 22
$ws.Range("H22").Value = 1636.8572
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1636.8572
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1636.8572
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = -2226.8572

# Row 27
$ws.Range("H27").Value = 1636.8572
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1636.8572
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1636.8572
$ws.Range("M27").Value = ""
$ws.Range("N27").Value = -1850.8572

# Row 40
$ws.Range("H40").Value = 4061.6667
$ws.Range("I40").Value = 4942.5713
$ws.Range("J40").Value = 2299.8572
$ws.Range("K40").Value = 4942.5713
$ws.Range("L40").Value = 2299.8572
$ws.Range("M40").Value = -4806.5713
$ws.Range("N40").Value = -2571.8572

# Row 122
$ws.Range("H122").Value = 8921.315000000001
$ws.Range("I122").Value = 10143
$ws.Range("J122").Value = 8357.462
$ws.Range("K122").Value = 30429
$ws.Range("L122").Value = 25072.386
$ws.Range("M122").Value = -27979
$ws.Range("N122").Value = -29972.386

# Row 126
$ws.Range("H126").Value = 7078.8
$ws.Range("I126").Value = 4916.625
$ws.Range("J126").Value = 9549.857
$ws.Range("K126").Value = 14749.875
$ws.Range("L126").Value = 28649.571
$ws.Range("M126").Value = -12279.875
$ws.Range("N126").Value = -33589.571

# Row 132
$ws.Range("H132").Value = 7468046
$ws.Range("I132").Value = 3749.843
$ws.Range("J132").Value = 31260490
$ws.Range("K132").Value = 11249.529
$ws.Range("L132").Value = 93781470
$ws.Range("M132").Value = -8719.528999999999
$ws.Range("N132").Value = -93786530

# Row 136
$ws.Range("H136").Value = 4706.7046
$ws.Range("I136").Value = 2190.4546
$ws.Range("J136").Value = 12255.454
$ws.Range("K136").Value = 6571.3638
$ws.Range("L136").Value = 36766.362
$ws.Range("M136").Value = -4021.3638
$ws.Range("N136").Value = -41866.362
